# Beta 0.9.1 - Expanded -selectr
#
# 1) Two label/header texts were corrected:
#      "Team-captain"       -> "Teamcaptain"
#      "Allergiën / dieet"  -> "Allergiën / Dieet"
#    These live in the header row (row 1) of Sheet1: column N ("Team-captain")
#    and column U ("Allergiën / dieet").
#
# 2) The visible selection/scroll position moved: the whole column X is now
#    selected (X1:X1048576, active cell X1) and the view is scrolled right so
#    column G becomes the left-most visible column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix the header text typos -----------------------------------------
$ws.Cells.Item(1, 14).Value = "Teamcaptain"        # N1: "Team-captain" -> "Teamcaptain"
$ws.Cells.Item(1, 21).Value = "Allergiën / Dieet"   # U1: "Allergiën / dieet" -> "Allergiën / Dieet"

# --- 2) Update the selection / scrolled view -------------------------------
# Select the entire column X (mirrors sqref="X1:X1048576", activeCell="X1")
$ws.Columns("X:X").Select()

# Scroll the window so column G is the left-most visible column
# (mirrors sheetView topLeftCell="G1")
$excel.ActiveWindow.ScrollColumn = 7
